# Auto-applied cell updates for cryptos sheet (D=Price text, E=Volume(1h) text)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''70.997.90'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +6.18%  '

$ws.Range("D3").Value = '''3.656.90'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +6.24%  '

$ws.Range("E4").Value = '  +0.06%  '

$ws.Range("D5").Value = '''596.44'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +2.96%  '

$ws.Range("D6").Value = '''194.70'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +3.35%  '

$ws.Range("D7").Value = '''0.647'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +2.59%  '

$ws.Range("D8").Value = '''3.649.95'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +6.19%  '

$ws.Range("E9").Value = '  +0.06%  '

$ws.Range("E10").Value = '  +8.01%  '

$ws.Range("E11").Value = '  +4.82%  '

$ws.Range("D12").Value = '''58.45'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +2.85%  '

$ws.Range("E13").Value = '  +6.59%  '

$ws.Range("D14").Value = '''9.98'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +5.81%  '

$ws.Range("D15").Value = '''4.243.49'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +6.41%  '

$ws.Range("D16").Value = '''20.13'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +7.42%  '

$ws.Range("D17").Value = '''3.660.05'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +6.45%  '

$ws.Range("D18").Value = '''71.021.93'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +6.35%  '

$ws.Range("D19").Value = '''12.83'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +6.31%  '

$ws.Range("E20").Value = '  +1.83%  '

$ws.Range("E21").Value = '  +4.31%  '

$ws.Range("D22").Value = '''489.43'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +1.24%  '

$ws.Range("D23").Value = '''19.00'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +12.26%  '

$ws.Range("E24").Value = '  -0.57%  '

$ws.Range("E25").Value = '  +4.30%  '

$ws.Range("D26").Value = '''91.51'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +2.54%  '

$ws.Range("E27").Value = '  +6.70%  '

$ws.Range("E28").Value = '  +4.46%  '

$ws.Range("D29").Value = '''9.63'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +6.59%  '

$ws.Range("D30").Value = '''32.91'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +5.43%  '

$ws.Range("D31").Value = '''7.82'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +5.97%  '

$ws.Range("D32").Value = '''0.123'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +9.77%  '

$ws.Range("D33").Value = '''629.84'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +5.57%  '

$ws.Range("E34").Value = '  +4.51%  '

$ws.Range("D35").Value = '''66.33'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +3.34%  '

$ws.Range("D36").Value = '''40.38'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +9.75%  '

$ws.Range("D37").Value = '''0.0₃0832'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +10.47%  '

$ws.Range("D38").Value = '''0.413'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +7.05%  '

$ws.Range("E39").Value = '  +0.30%  '

$ws.Range("E40").Value = '  +0.03%  '

$ws.Range("E41").Value = '  +2.35%  '

$ws.Range("D42").Value = '''3.329.50'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +3.09%  '

$ws.Range("E43").Value = '  +9.51%  '

$ws.Range("D44").Value = '''2.82'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +11.36%  '

$ws.Range("D45").Value = '''3.07'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +8.29%  '

$ws.Range("D46").Value = '''0.0456'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +6.09%  '

$ws.Range("D47").Value = '''9.47'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +10.20%  '

$ws.Range("E48").Value = '  +3.98%  '

$ws.Range("D49").Value = '''3.32'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +2.90%  '

$ws.Range("D50").Value = '''3.27'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -2.98%  '

$ws.Range("D51").Value = '''1.00'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +0.08%  '
